$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1541.0588
$ws.Range("I112").Value = 1700
$ws.Range("J112").Value = 1531.125
$ws.Range("K112").Value = 5100
$ws.Range("L112").Value = 4593.375
$ws.Range("M112").Value = -3992
$ws.Range("N112").Value = -6809.375
$ws.Range("H129").Value = 1248.5
$ws.Range("I129").Value = 997
$ws.Range("J129").Value = 1500
$ws.Range("K129").Value = 2991
$ws.Range("L129").Value = 4500
$ws.Range("M129").Value = 2009
$ws.Range("N129").Value = -14500
$ws.Range("H137").Value = 1131.8877
$ws.Range("I137").Value = 717.63635
$ws.Range("J137").Value = 1376
$ws.Range("K137").Value = 2152.90905
$ws.Range("L137").Value = 4128
$ws.Range("M137").Value = 397.0909499999998
$ws.Range("N137").Value = -9228
$ws.Range("H138").Value = 7039.359
$ws.Range("I138").Value = 2124.25
$ws.Range("J138").Value = 7601.086
$ws.Range("K138").Value = 6372.75
$ws.Range("L138").Value = 22803.258
$ws.Range("M138").Value = -1232.75
$ws.Range("N138").Value = -33083.258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6546.9487
$ws.Range("I32").Value = 3897.4062
$ws.Range("J32").Value = 18659.143
$ws.Range("K32").Value = 3897.4062
$ws.Range("L32").Value = 18659.143
$ws.Range("M32").Value = -3610.4062
$ws.Range("N32").Value = -19233.143
$ws.Range("H61").Value = 1652.6364
$ws.Range("I61").Value = 1518.0526
$ws.Range("J61").Value = 2505
$ws.Range("K61").Value = 1518.0526
$ws.Range("L61").Value = 2505
$ws.Range("M61").Value = -1306.0526
$ws.Range("N61").Value = -2929
$ws.Range("H74").Value = 1002.64105
$ws.Range("I74").Value = 519.2105
$ws.Range("J74").Value = 1461.9
$ws.Range("K74").Value = 519.2105
$ws.Range("L74").Value = 1461.9
$ws.Range("M74").Value = 354.7895
$ws.Range("N74").Value = -3209.9
$ws.Range("H77").Value = 1002.64105
$ws.Range("I77").Value = 519.2105
$ws.Range("J77").Value = 1461.9
$ws.Range("K77").Value = 2596.0525
$ws.Range("L77").Value = 7309.5
$ws.Range("M77").Value = 1771.9475
$ws.Range("N77").Value = -16045.5
$ws.Range("H109").Value = 20377
$ws.Range("J109").Value = 20377
$ws.Range("L109").Value = 20377
$ws.Range("N109").Value = -23151
$ws.Range("H132").Value = 2697.5945
$ws.Range("I132").Value = 1783.8334
$ws.Range("J132").Value = 3563.2632
$ws.Range("K132").Value = 5351.5002
$ws.Range("L132").Value = 10689.7896
$ws.Range("M132").Value = -2821.5002
$ws.Range("N132").Value = -15749.7896
$ws.Range("H136").Value = 1652.6364
$ws.Range("I136").Value = 1518.0526
$ws.Range("J136").Value = 2505
$ws.Range("K136").Value = 4554.1578
$ws.Range("L136").Value = 7515
$ws.Range("M136").Value = -2004.1578
$ws.Range("N136").Value = -12615

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 10040
$ws.Range("I113").Value = 10040
$ws.Range("K113").Value = 10040
$ws.Range("M113").Value = -7870
$ws.Range("H134").Value = 1915.0667
$ws.Range("I134").Value = 1790.0769
$ws.Range("J134").Value = 2727.5
$ws.Range("K134").Value = 5370.2307
$ws.Range("L134").Value = 8182.5
$ws.Range("M134").Value = -2835.2307
$ws.Range("N134").Value = -13252.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1295.4154
$ws.Range("I31").Value = 829.1429000000001
$ws.Range("J31").Value = 1423.4117
$ws.Range("K31").Value = 829.1429000000001
$ws.Range("L31").Value = 1423.4117
$ws.Range("M31").Value = -534.1429000000001
$ws.Range("N31").Value = -2013.4117
$ws.Range("H34").Value = 1295.4154
$ws.Range("I34").Value = 829.1429000000001
$ws.Range("J34").Value = 1423.4117
$ws.Range("K34").Value = 829.1429000000001
$ws.Range("L34").Value = 1423.4117
$ws.Range("M34").Value = -627.1429000000001
$ws.Range("N34").Value = -1827.4117
$ws.Range("H132").Value = 4585.1113
$ws.Range("I132").Value = 4214
$ws.Range("K132").Value = 12642
$ws.Range("M132").Value = -10112
$ws.Range("H134").Value = 2261.3
$ws.Range("I134").Value = 2039.3846
$ws.Range("J134").Value = 2673.4285
$ws.Range("K134").Value = 6118.1538
$ws.Range("L134").Value = 8020.2855
$ws.Range("M134").Value = -3583.1538
$ws.Range("N134").Value = -13090.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7886.9414
$ws.Range("I132").Value = 5151.433
$ws.Range("J132").Value = 28403.25
$ws.Range("K132").Value = 15454.299
$ws.Range("L132").Value = 85209.75
$ws.Range("M132").Value = -12924.299
$ws.Range("N132").Value = -90269.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 10000
$ws.Range("J11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("N11").Value = -10280
$ws.Range("H14").Value = 3000
$ws.Range("I14").Value = 3000
$ws.Range("K14").Value = 3000
$ws.Range("M14").Value = -2828
$ws.Range("H40").Value = 2208.3704
$ws.Range("I40").Value = 2124.8096
$ws.Range("J40").Value = 2500.8333
$ws.Range("K40").Value = 2124.8096
$ws.Range("L40").Value = 2500.8333
$ws.Range("M40").Value = -1988.8096
$ws.Range("N40").Value = -2772.8333
$ws.Range("H46").Value = 1011.44446
$ws.Range("I46").Value = 1004.86365
$ws.Range("J46").Value = 1040.4
$ws.Range("K46").Value = 1004.86365
$ws.Range("L46").Value = 1040.4
$ws.Range("M46").Value = -816.86365
$ws.Range("N46").Value = -1416.4
$ws.Range("H55").Value = 181.78261
$ws.Range("I55").Value = 140.5625
$ws.Range("J55").Value = 276
$ws.Range("K55").Value = 140.5625
$ws.Range("L55").Value = 276
$ws.Range("M55").Value = 32.4375
$ws.Range("N55").Value = -622
$ws.Range("H132").Value = 2975.4443
$ws.Range("I132").Value = 2436.3333
$ws.Range("J132").Value = 3649.3333
$ws.Range("K132").Value = 7308.999899999999
$ws.Range("L132").Value = 10947.9999
$ws.Range("M132").Value = -4778.999899999999
$ws.Range("N132").Value = -16007.9999
$ws.Range("H136").Value = 2281.238
$ws.Range("I136").Value = 1567
$ws.Range("J136").Value = 2816.9167
$ws.Range("K136").Value = 4701
$ws.Range("L136").Value = 8450.750100000001
$ws.Range("M136").Value = -2151
$ws.Range("N136").Value = -13550.7501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1753.6061
$ws.Range("I132").Value = 1176.4445
$ws.Range("J132").Value = 2446.2
$ws.Range("K132").Value = 3529.3335
$ws.Range("L132").Value = 7338.599999999999
$ws.Range("M132").Value = -999.3335000000002
$ws.Range("N132").Value = -12398.6
